$p = $ppt.ActivePresentation

# The author opened the Notes pane on the title slide (slide 1) during
# class and left it blank -- this provisions a (currently empty) Notes
# Page / notes slide part for that slide, without touching any of the
# slide's own visible content or the notes already present on the later
# "Dataset 3" slides.
$s1 = $p.Slides.Item(1)
$notesPage = $s1.NotesPage
$notesBody = $notesPage.Shapes.AddPlaceholder(2)
